# The data row 2 is turned into a "Total" row: the value "5000" previously
# held in M2 (valorDescarga) is replaced with the label "Total".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "Total"

# The saved cursor/selection moves from M2 down to M3.
[void]$ws.Range("M3").Select()
